$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.589.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.842.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.46%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4242"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.56%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3641"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.26"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07258"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8931"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.64"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.840.42"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.571"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.355"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06875"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "78.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.06%  "
$ws.Range("E19").Value = "  -2.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("E21").Value = "  -2.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.583.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.984"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.40%  "
$ws.Range("E24").Value = "  -2.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.044.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.48%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.021"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.229"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.835"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08908"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7805"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.48%  "
$ws.Range("E34").Value = "  -5.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.962"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.84%  "
$ws.Range("E36").Value = "  -6.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9995"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05406"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.097"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01922"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.780"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.845"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5067"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.47%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1649"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.223"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06613"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.79%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4707"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.23%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9985"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.628"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.51%  "
